# Update crypto price/volume figures (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="27.025.65"; E="  -3.31%  " }
    @{ Row=3; D="1.715.32"; E="  -3.04%  " }
    @{ Row=4; D="1.001"; E="  +0.05%  " }
    @{ Row=5; D="309.32"; E="  -5.89%  " }
    @{ Row=6; D="1.001"; E="  +0.09%  " }
    @{ Row=7; D="0.4753"; E="  +4.40%  " }
    @{ Row=8; D="0.3466"; E="  -1.92%  " }
    @{ Row=9; D="41.93"; E="  -0.26%  " }
    @{ Row=10; D="0.07238"; E="  -2.04%  " }
    @{ Row=11; D="1.039"; E="  -5.26%  " }
    @{ Row=12; D="1.001"; E="  +0.09%  " }
    @{ Row=13; D="19.77"; E="  -4.82%  " }
    @{ Row=14; D="5.828"; E="  -3.10%  " }
    @{ Row=15; D="1.722.59"; E="  -2.93%  " }
    @{ Row=16; D="6.835"; E="  -4.97%  " }
    @{ Row=17; D="86.85"; E="  -6.26%  " }
    @{ Row=18; D=$null; E="  -2.43%  " }
    @{ Row=19; D="0.06370"; E="  -1.11%  " }
    @{ Row=20; D=$null; E="  +0.07%  " }
    @{ Row=21; D="16.47"; E="  -2.91%  " }
    @{ Row=22; D="5.612"; E="  -2.83%  " }
    @{ Row=23; D="27.075.40"; E="  -3.20%  " }
    @{ Row=24; D="10.73"; E="  -4.37%  " }
    @{ Row=25; D=$null; E="  -0.11%  " }
    @{ Row=26; D="19.93"; E="  -1.02%  " }
    @{ Row=27; D="150.80"; E="  -5.47%  " }
    @{ Row=28; D="1.911.14"; E="  -3.38%  " }
    @{ Row=29; D="2.063"; E="  -4.29%  " }
    @{ Row=30; D="120.35"; E="  -3.14%  " }
    @{ Row=31; D="1.025"; E="  -5.03%  " }
    @{ Row=32; D="0.09129"; E="  -0.90%  " }
    @{ Row=33; D="3.601"; E="  -1.60%  " }
    @{ Row=34; D="5.301"; E="  -5.57%  " }
    @{ Row=35; D="1.466"; E="  +6.21%  " }
    @{ Row=36; D="0.02174"; E="  -4.84%  " }
    @{ Row=37; D="0.05843"; E="  -4.56%  " }
    @{ Row=38; D="0.1997"; E="  -4.58%  " }
    @{ Row=39; D="10.91"; E="  -8.07%  " }
    @{ Row=40; D="1.000"; E="  +0.00%  " }
    @{ Row=41; D="4.702"; E="  -5.01%  " }
    @{ Row=42; D="0.5961"; E="  -4.76%  " }
    @{ Row=43; D="1.080"; E="  -8.28%  " }
    @{ Row=44; D="7.473"; E="  -4.55%  " }
    @{ Row=45; D="12.68"; E="  -4.53%  " }
    @{ Row=46; D="3.577"; E="  -4.23%  " }
    @{ Row=47; D="0.5563"; E="  -4.92%  " }
    @{ Row=48; D="118.83"; E="  -3.14%  " }
    @{ Row=49; D="1.822"; E="  -5.92%  " }
    @{ Row=50; D="1.115"; E="  -1.46%  " }
    @{ Row=51; D="0.06627"; E="  -2.91%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
}
